$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H2 from 7 to 6.3 (minimum 101/110L grade)
$ws.Range("H2").Value = 6.3

# Update the active selection to H2 (matches saved selection in the file)
$ws.Range("H2").Select()
